$d = $word.ActiveDocument

# Append a new bulleted paragraph analysing the breach against the CIA triad.
$p1 = $d.Paragraphs.Add()
$p1.Range.Text = "Analysis: The breach complied with the confidentiality pillar of the CIA triad, as the data remained encrypted and confidential, however it was encrypted by the hackers, and so the data was hidden to the victims, not to the hackers. This meant that the availability pillar was not followed as the data was not available for the victim to use. The integrity of the data was also compromised as the victim can not be certain the data was not tampered with as they do not have access to it."

# Append a further bulleted paragraph with a recommendation for victims.
$p2 = $d.Paragraphs.Add()
$p2.Range.Text = "I would recommend that a victim do some research online to see if there was a cure to the virus available online, however if it was not available, and if it was the only way for life to continue and less money to be wasted, I would recommend paying for the payload so that life can continue as quickly as possible."

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
